# Auto-generated edit script: applies scheduled market-data refresh values
# to the Leve profit calculator sheets (columns H-N) per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2608.5
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H47").Value = 5555
$ws.Range("I47").Value = 5555
$ws.Range("K47").Value = 5555
$ws.Range("M47").Value = -4583
$ws.Range("H137").Value = 8896.071
$ws.Range("I137").Value = 1844.7
$ws.Range("K137").Value = 5534.1
$ws.Range("M137").Value = -2984.1
$ws.Range("H138").Value = 24090.4
$ws.Range("I138").Value = 49668.617
$ws.Range("J138").Value = 5568.241
$ws.Range("K138").Value = 149005.851
$ws.Range("L138").Value = 16704.723
$ws.Range("M138").Value = -143865.851
$ws.Range("N138").Value = -26984.723

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3342.877
$ws.Range("I32").Value = 3134.7302
$ws.Range("K32").Value = 3134.7302
$ws.Range("M32").Value = -2847.7302
$ws.Range("H74").Value = 50972.332
$ws.Range("I74").Value = 92547.17999999999
$ws.Range("K74").Value = 92547.17999999999
$ws.Range("M74").Value = -91673.17999999999
$ws.Range("H77").Value = 50972.332
$ws.Range("I77").Value = 92547.17999999999
$ws.Range("K77").Value = 462735.9
$ws.Range("M77").Value = -458367.9
$ws.Range("H122").Value = 2541.5557
$ws.Range("I122").Value = 2482.5
$ws.Range("K122").Value = 7447.5
$ws.Range("M122").Value = -4997.5
$ws.Range("H132").Value = 3581.2092
$ws.Range("I132").Value = 3469.5588
$ws.Range("K132").Value = 10408.6764
$ws.Range("M132").Value = -7878.6764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H105").Value = 2449.8572
$ws.Range("I105").Value = 2089.8
$ws.Range("K105").Value = 2089.8
$ws.Range("M105").Value = -342.8000000000002
$ws.Range("H107").Value = 2039.7667
$ws.Range("I107").Value = 1653.875
$ws.Range("J107").Value = 3583.3333
$ws.Range("K107").Value = 1653.875
$ws.Range("L107").Value = 3583.3333
$ws.Range("M107").Value = 266.125
$ws.Range("N107").Value = -7423.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1475
$ws.Range("J25").Value = 1475
$ws.Range("L25").Value = 1475
$ws.Range("N25").Value = -1823
$ws.Range("H31").Value = 33442.97
$ws.Range("I31").Value = 58030.445
$ws.Range("K31").Value = 58030.445
$ws.Range("M31").Value = -57735.445
$ws.Range("H34").Value = 33442.97
$ws.Range("I34").Value = 58030.445
$ws.Range("K34").Value = 58030.445
$ws.Range("M34").Value = -57828.445
$ws.Range("H58").Value = 3442.825
$ws.Range("I58").Value = 3051.4138
$ws.Range("K58").Value = 3051.4138
$ws.Range("M58").Value = -2848.4138
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H125").Value = 109833.336
$ws.Range("J125").Value = 109833.336
$ws.Range("L125").Value = 109833.336
$ws.Range("N125").Value = -114753.336
$ws.Range("H132").Value = 173055.64
$ws.Range("I132").Value = 212361.55
$ws.Range("K132").Value = 637084.6499999999
$ws.Range("M132").Value = -634554.6499999999
$ws.Range("H134").Value = 26805.846
$ws.Range("I134").Value = 22657.725
$ws.Range("J134").Value = 38835.4
$ws.Range("K134").Value = 67973.17499999999
$ws.Range("L134").Value = 116506.2
$ws.Range("M134").Value = -65438.17499999999
$ws.Range("N134").Value = -121576.2
$ws.Range("H136").Value = 3442.825
$ws.Range("I136").Value = 3051.4138
$ws.Range("K136").Value = 9154.241399999999
$ws.Range("M136").Value = -6604.241399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1379.8
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 1966.3334
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 5899.0002
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = -6123.0002
$ws.Range("H68").Value = 12501500
$ws.Range("J68").Value = 50000000
$ws.Range("L68").Value = 150000000
$ws.Range("N68").Value = -150001622
$ws.Range("H71").Value = 12501500
$ws.Range("J71").Value = 50000000
$ws.Range("L71").Value = 450000000
$ws.Range("N71").Value = -450008112
$ws.Range("H92").Value = 2500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2500
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -9996
$ws.Range("H131").Value = 35786.332
$ws.Range("I131").Value = 144042.14
$ws.Range("K131").Value = 432126.42
$ws.Range("M131").Value = -427086.42
$ws.Range("H135").Value = 1379.8
$ws.Range("I135").Value = 500
$ws.Range("J135").Value = 1966.3334
$ws.Range("K135").Value = 4500
$ws.Range("L135").Value = 17697.0006
$ws.Range("M135").Value = -1965
$ws.Range("N135").Value = -22767.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H126").Value = 21813.715
$ws.Range("I126").Value = 35056.855
$ws.Range("K126").Value = 105170.565
$ws.Range("M126").Value = -102700.565
$ws.Range("H132").Value = 4352.9688
$ws.Range("I132").Value = 4299.9287
$ws.Range("J132").Value = 4724.25
$ws.Range("K132").Value = 12899.7861
$ws.Range("L132").Value = 14172.75
$ws.Range("M132").Value = -10369.7861
$ws.Range("N132").Value = -19232.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9173.429
$ws.Range("I7").Value = 9843
$ws.Range("K7").Value = 9843
$ws.Range("M7").Value = -9731
$ws.Range("H22").Value = 3522.3333
$ws.Range("I22").Value = 2450.1667
$ws.Range("J22").Value = 5666.6665
$ws.Range("K22").Value = 2450.1667
$ws.Range("L22").Value = 5666.6665
$ws.Range("M22").Value = -2155.1667
$ws.Range("N22").Value = -6256.6665
$ws.Range("H27").Value = 3522.3333
$ws.Range("I27").Value = 2450.1667
$ws.Range("J27").Value = 5666.6665
$ws.Range("K27").Value = 2450.1667
$ws.Range("L27").Value = 5666.6665
$ws.Range("M27").Value = -2343.1667
$ws.Range("N27").Value = -5880.6665
$ws.Range("H40").Value = 3690.394
$ws.Range("I40").Value = 3285.7036
$ws.Range("J40").Value = 5511.5
$ws.Range("K40").Value = 3285.7036
$ws.Range("L40").Value = 5511.5
$ws.Range("M40").Value = -3149.7036
$ws.Range("N40").Value = -5783.5
$ws.Range("H43").Value = 30000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H60").Value = 71059.5
$ws.Range("J60").Value = 71059.5
$ws.Range("L60").Value = 71059.5
$ws.Range("N60").Value = -72077.5
$ws.Range("H126").Value = 9173.429
$ws.Range("I126").Value = 9843
$ws.Range("K126").Value = 29529
$ws.Range("M126").Value = -27059
$ws.Range("H136").Value = 46413.39
$ws.Range("I136").Value = 58450.777
$ws.Range("K136").Value = 175352.331
$ws.Range("M136").Value = -172802.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 150527.78
$ws.Range("I2").Value = 178964.28
$ws.Range("J2").Value = 51000
$ws.Range("K2").Value = 178964.28
$ws.Range("L2").Value = 51000
$ws.Range("M2").Value = -178852.28
$ws.Range("N2").Value = -51224
$ws.Range("H14").Value = 779.6667
$ws.Range("I14").Value = 649.4
$ws.Range("J14").Value = 1186.75
$ws.Range("K14").Value = 649.4
$ws.Range("L14").Value = 1186.75
$ws.Range("M14").Value = -481.4
$ws.Range("N14").Value = -1522.75
$ws.Range("H58").Value = 44273.832
$ws.Range("I58").Value = 42000
$ws.Range("J58").Value = 44728.6
$ws.Range("K58").Value = 42000
$ws.Range("L58").Value = 44728.6
$ws.Range("M58").Value = -41692
$ws.Range("N58").Value = -45344.6
$ws.Range("H105").Value = 36403.25
$ws.Range("J105").Value = 36403.25
$ws.Range("L105").Value = 36403.25
$ws.Range("N105").Value = -43391.25
$ws.Range("H122").Value = 3999
$ws.Range("I122").Value = 3598.8
$ws.Range("K122").Value = 10796.4
$ws.Range("M122").Value = -8346.400000000001
$ws.Range("H126").Value = 5710.8335
$ws.Range("I126").Value = 7224.353
$ws.Range("J126").Value = 2035.1428
$ws.Range("K126").Value = 21673.059
$ws.Range("L126").Value = 6105.428400000001
$ws.Range("M126").Value = -19203.059
$ws.Range("N126").Value = -11045.4284
$ws.Range("H132").Value = 2222.5
$ws.Range("I132").Value = 2171.75
$ws.Range("K132").Value = 6515.25
$ws.Range("M132").Value = -3985.25

